# Updated symbol list on Tue Dec 27 19:26:42 UTC 2022 with GitHub Actions
#
# This script re-applies the scraped price/rank refresh to the cryptos
# worksheet. Column D values are stored as text (to preserve literal
# formatting such as trailing zeros), so numeric-looking strings are
# written with a leading quote-prefix and the cell style is reset back
# to "Normal" afterwards so no stray number-format/style is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [object]$Worksheet,
        [string]$CellRef,
        [string]$Text
    )
    $cell = $Worksheet.Range($CellRef)
    $cell.Value = "'" + $Text
    $cell.Style = "Normal"
}

# ---- Simple price refreshes (column D only) ----
Set-TextValue $ws "D2"  "245.75"
Set-TextValue $ws "D3"  "23.75"
Set-TextValue $ws "D4"  "5.335"
Set-TextValue $ws "D7"  "3.350"
Set-TextValue $ws "D8"  "0.8113"
Set-TextValue $ws "D9"  "0.9207"
Set-TextValue $ws "D11" "0.07349"
Set-TextValue $ws "D12" "0.03066"
Set-TextValue $ws "D13" "0.03075"
Set-TextValue $ws "D14" "0.09362"
Set-TextValue $ws "D15" "3.854"
Set-TextValue $ws "D16" "0.001550"
Set-TextValue $ws "D17" "0.04677"

# ---- Rank shuffle rows 18-24: "One" jumps to rank 18, others shift down ----
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D18" "0.0005998"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D19" "0.006063"
$ws.Range("E19").Value = "18TigerCashTCH"

$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws "D20" "0.001244"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws "D21" "0.004688"
$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws "D22" "0.00008807"
$ws.Range("E22").Value = "21NitroExNTXBestin24h"

$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D23" "3.593"
$ws.Range("E23").Value = "22LEOLEO"

$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws "D24" "2.158"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# ---- More simple price refreshes ----
Set-TextValue $ws "D25" "0.3229"
Set-TextValue $ws "D40" "0.03837"

# ---- Rank shuffle rows 41-43: KickToken moves to rank 43 (worst in 24h) ----
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D41" "0.1066"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D42" "0.002702"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D43" "0.003096"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# ---- Final simple price refreshes ----
Set-TextValue $ws "D44" "0.007780"
Set-TextValue $ws "D45" "0.00005256"
Set-TextValue $ws "D47" "0.6810"
Set-TextValue $ws "D48" "0.001859"
$ws.Range("E48").Value = "47BOLOBOLO"
Set-TextValue $ws "D49" "0.00002103"
